$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: corrected Hydrogen value (Iron & steel column)
$ws.Range("B3").Value = 4587760.284859059

# D3: Hydrogen / Non-metallic minerals value removed (corrected to blank)
$ws.Range("D3").Value = ""

# C4: corrected Methanol value (Chemicals column)
$ws.Range("C4").Value = 193.9960356878017

# C5: corrected Ammonia value (Chemicals column)
$ws.Range("C5").Value = 5592.422905162121

# Row 7 label was "Other", corrected code renames it to "Biogas"
# and updates its Non-metallic minerals value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 5334.58045222516

# A new row 8 is inserted for "Other", carrying the remaining
# Non-metallic minerals value previously lumped into row 7
$ws.Range("A8").Value = "Other"

# Copy the row-label formatting (bold, centered, bordered) from A7 to A8
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D8").Value = 4181.743697986703
